# Apply crypto price/volume updates (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to be written/stored as literal text (matches the
    # source data, which is all inlineStr/shared-string, never numeric),
    # then drop the temporary text number-format so the cell keeps the
    # workbook default style (no explicit s= attribute), same as before.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "63.522.91"
$ws.Range("E2").Value = "  -3.60%  "
Set-TextValue $ws.Range("D3") "2.604.47"
$ws.Range("E3").Value = "  -2.38%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "571.51"
$ws.Range("E5").Value = "  -4.55%  "
Set-TextValue $ws.Range("D6") "154.56"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E7").Value = "  +0.00%  "
Set-TextValue $ws.Range("D8") "0.619"
$ws.Range("E8").Value = "  -5.30%  "
$ws.Range("E9").Value = "  -7.26%  "
$ws.Range("E10").Value = "  -0.73%  "
Set-TextValue $ws.Range("D11") "0.381"
$ws.Range("E11").Value = "  -5.32%  "
$ws.Range("E12").Value = "  -0.48%  "
Set-TextValue $ws.Range("D13") "28.16"
$ws.Range("E13").Value = "  -3.10%  "
Set-TextValue $ws.Range("D14") "3.076.78"
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("E15").Value = "  -8.89%  "
Set-TextValue $ws.Range("D16") "63.349.25"
$ws.Range("E16").Value = "  -3.71%  "
Set-TextValue $ws.Range("D17") "2.657.09"
$ws.Range("E17").Value = "  +3.46%  "
Set-TextValue $ws.Range("D18") "11.94"
$ws.Range("E18").Value = "  -5.38%  "
$ws.Range("E19").Value = "  +0.46%  "
Set-TextValue $ws.Range("D20") "4.52"
$ws.Range("E20").Value = "  -5.90%  "
Set-TextValue $ws.Range("D21") "342.45"
$ws.Range("E21").Value = "  -2.32%  "
$ws.Range("E22").Value = "  +0.00%  "
Set-TextValue $ws.Range("D23") "67.01"
$ws.Range("E23").Value = "  -4.00%  "
Set-TextValue $ws.Range("D24") "1.77"
$ws.Range("E24").Value = "  -3.25%  "
$ws.Range("E25").Value = "  -4.42%  "
Set-TextValue $ws.Range("D26") "586.24"
$ws.Range("E26").Value = "  +3.32%  "
Set-TextValue $ws.Range("D27") "9.09"
$ws.Range("E27").Value = "  -5.45%  "
$ws.Range("E28").Value = "  -4.43%  "
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").Value = "  -2.29%  "
Set-TextValue $ws.Range("D31") "7.88"
$ws.Range("E31").Value = "  -4.18%  "
$ws.Range("E32").Value = "  -4.37%  "
Set-TextValue $ws.Range("D33") "1.72"
$ws.Range("E33").Value = "  -5.80%  "
Set-TextValue $ws.Range("D34") "6.50"
$ws.Range("E34").Value = "  -3.03%  "
$ws.Range("E35").Value = "  -2.97%  "
Set-TextValue $ws.Range("D36") "0.403"
$ws.Range("E36").Value = "  -4.75%  "
Set-TextValue $ws.Range("D37") "0.999"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  -4.65%  "
Set-TextValue $ws.Range("D39") "154.81"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("E40").Value = "  -4.93%  "
Set-TextValue $ws.Range("D41") "1.00"
$ws.Range("E41").Value = "  +0.00%  "
Set-TextValue $ws.Range("D42") "41.28"
$ws.Range("E42").Value = "  -3.32%  "
Set-TextValue $ws.Range("D43") "2.43"
$ws.Range("E43").Value = "  +4.60%  "
Set-TextValue $ws.Range("D44") "155.69"
$ws.Range("E44").Value = "  -3.16%  "
Set-TextValue $ws.Range("D45") "3.89"
$ws.Range("E45").Value = "  -5.26%  "
Set-TextValue $ws.Range("D46") "22.99"
$ws.Range("E46").Value = "  +0.12%  "
Set-TextValue $ws.Range("D47") "0.0585"
$ws.Range("E47").Value = "  -5.62%  "
$ws.Range("E48").Value = "  -2.45%  "
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("E50").Value = "  -3.98%  "
Set-TextValue $ws.Range("D51") "18.65"
$ws.Range("E51").Value = "  -6.27%  "
